$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Astro solar-production model was retrained; the forecast horizon for this
# file has moved forward by 2 days. Shift every timestamp in column A (rows
# 2-97) forward by 2 days accordingly. Target values are supplied explicitly
# (rather than read-modify-write) to avoid any COM round-trip precision loss.
$newTimestamps = @(46046.01041666666, 46046.02083333334, 46046.03125, 46046.04166666666, 46046.05208333334, 46046.0625, 46046.07291666666, 46046.08333333334, 46046.09375, 46046.10416666666, 46046.11458333334, 46046.125, 46046.13541666666, 46046.14583333334, 46046.15625, 46046.16666666666, 46046.17708333334, 46046.1875, 46046.19791666666, 46046.20833333334, 46046.21875, 46046.22916666666, 46046.23958333334, 46046.25, 46046.26041666666, 46046.27083333334, 46046.28125, 46046.29166666666, 46046.30208333334, 46046.3125, 46046.32291666666, 46046.33333333334, 46046.34375, 46046.35416666666, 46046.36458333334, 46046.375, 46046.38541666666, 46046.39583333334, 46046.40625, 46046.41666666666, 46046.42708333334, 46046.4375, 46046.44791666666, 46046.45833333334, 46046.46875, 46046.47916666666, 46046.48958333334, 46046.5, 46046.51041666666, 46046.52083333334, 46046.53125, 46046.54166666666, 46046.55208333334, 46046.5625, 46046.57291666666, 46046.58333333334, 46046.59375, 46046.60416666666, 46046.61458333334, 46046.625, 46046.63541666666, 46046.64583333334, 46046.65625, 46046.66666666666, 46046.67708333334, 46046.6875, 46046.69791666666, 46046.70833333334, 46046.71875, 46046.72916666666, 46046.73958333334, 46046.75, 46046.76041666666, 46046.77083333334, 46046.78125, 46046.79166666666, 46046.80208333334, 46046.8125, 46046.82291666666, 46046.83333333334, 46046.84375, 46046.85416666666, 46046.86458333334, 46046.875, 46046.88541666666, 46046.89583333334, 46046.90625, 46046.91666666666, 46046.92708333334, 46046.9375, 46046.94791666666, 46046.95833333334, 46046.96875, 46046.97916666666, 46046.98958333334, 46047.0)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value2 = $newTimestamps[$i]
}

# The production reading for the 07:15-07:30 interval (row 31) was an
# erroneous spike from the old model; the retrained model now predicts 0.
$ws.Range("B31").Value2 = 0
